# "warn when closing with unsaved data" (Id 12) is done: move it from the
# Active (Todo) sheet to the Inactive (Done) sheet as the newest completed
# item, with a Done date of 3/5/2018.

$wb = $excel.ActiveWorkbook
$active = $wb.Worksheets.Item("Active")
$inactive = $wb.Worksheets.Item("Inactive")

# Remove the task row from the Active sheet (row 3: Id 12, "warn when
# closing with unsaved data").
$active.Rows("3:3").Delete()

# Insert a new row right under the header of the Inactive sheet and fill it
# in as the newly completed task. Clear the formatting it inherits from the
# header row so it matches the plain data rows below it.
$inactive.Rows("2:2").Insert()
$inactive.Rows("2:2").ClearFormats()

$inactive.Cells.Item(2, 1).Value = 12
$inactive.Cells.Item(2, 2).Value = "warn when closing with unsaved data"
$inactive.Cells.Item(2, 3).Value = "Done"
$inactive.Cells.Item(2, 4).Value = "Bug"

# Created/Done columns hold plain text dates (e.g. "12/8/2017") in this
# workbook, not real Excel date values - force text formatting first so
# Excel doesn't silently convert them into date serial numbers, then clear
# the format back off so the cell matches its plain, unstyled neighbours.
$inactive.Cells.Item(2, 5).NumberFormat = "@"
$inactive.Cells.Item(2, 5).Value = "12/8/2017"
$inactive.Cells.Item(2, 5).ClearFormats()

$inactive.Cells.Item(2, 6).NumberFormat = "@"
$inactive.Cells.Item(2, 6).Value = "3/5/2018"
$inactive.Cells.Item(2, 6).ClearFormats()
